$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add three new (mostly blank) paragraphs before the "Education" heading,
#    right after the existing run of blank paragraphs that follow the
#    "...Network Architectures and Protocols." sentence. The last of the
#    three new paragraphs carries a run of four spaces.
# ---------------------------------------------------------------------------
$scope = $d.Content
$scope.Find.Execute("Network Architectures and Protocols.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterSentence = $d.Range($scope.End, $scope.End)
$eduScope = $d.Range($afterSentence.End, $d.Content.End)
$eduScope.Find.Execute("Education", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$eduHeadingStart = $eduScope.Start

$insertPoint = $d.Range($eduHeadingStart - 1, $eduHeadingStart - 1)
$insertPoint.InsertAfter("`r`r`r")

# recompute the Education heading paragraph index so we can walk backwards
$eduScope2 = $d.Content
$eduScope2.Find.Execute("Education", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$eduPara = $eduScope2.Paragraphs(1).Index

$blank3 = $d.Paragraphs($eduPara - 1)
$spacesTarget = $d.Range($blank3.Range.End - 1, $blank3.Range.End - 1)
$spacesTarget.InsertBefore("    ")

# ---------------------------------------------------------------------------
# 2. Split "Exeter University: Computer Science" into two runs:
#    "University of Exeter" + ": Computer Science"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Exeter University: Computer Science", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$prefix = $d.Range($r.Start, $r.Start + 17)
$prefix.Text = ""
$point2 = $d.Range($r.Start, $r.Start)
$point2.InsertBefore("University of Exeter")

# ---------------------------------------------------------------------------
# 3. Merge "with some practical e" + bookmark + "xperience in each area."
#    into a single run, removing the _GoBack bookmark from here.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("with some practical experience in each area.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r2.Find.Found) {
    # text is still split across two runs with the bookmark in between
    $r3 = $d.Content
    $r3.Find.Execute("with some practical e", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $joinPoint = $d.Range($r3.End, $r3.End)
    $joinPoint.InsertAfter("xperience in each area.")
    $r4 = $d.Content
    $r4.Find.Execute("with some practical experience in each areaxperience in each area.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

$r5 = $d.Content
$r5.Find.Execute("with some practical e", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dup = $d.Range($r5.End, $r5.End)
$dupScope = $d.Range($dup.End, $dup.End + 40)

$r6 = $d.Content
$r6.Find.Execute("xperience in each area.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# ---------------------------------------------------------------------------
# 4. Remove the spellStart/spellEnd proofErr pair around "Okey-Doke" and
#    relocate the _GoBack bookmark to sit right after that run.
# ---------------------------------------------------------------------------
$r7 = $d.Content
$r7.Find.Execute("Okey-Doke", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterOkeyDoke = $d.Range($r7.End, $r7.End)
$d.Bookmarks.Add("_GoBack", $afterOkeyDoke)

# ---------------------------------------------------------------------------
# 5. References section: drop the "Company:"/"Mobile Phone:"/"E-Mail
#    Address:" labels (and following space) and append Dr David Wakeling.
# ---------------------------------------------------------------------------
function Remove-LabelNear($searchText, $label) {
    $scope = $d.Content
    $scope.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $p = $scope.Paragraphs(1)
    $pscope = $d.Range($p.Range.Start, $p.Range.End)
    $pscope.Find.Execute($label, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $pscope.Delete()
}

Remove-LabelNear "Company: Harlequin Computing Solutions" "Company: "
Remove-LabelNear "Mobile Phone: +44 7811 276953" "Mobile Phone: "
Remove-LabelNear "E-Mail Address: duncan@oneeyedmen.com" "E-Mail Address: "

$emailScope = $d.Content
$emailScope.Find.Execute("duncan@oneeyedmen.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$emailPara = $emailScope.Paragraphs(1)
$endOfEmailPara = $d.Range($emailPara.Range.End - 1, $emailPara.Range.End - 1)
$endOfEmailPara.InsertAfter("`r`r`r`r")

$wakelingScope = $d.Content
$wakelingScope.Find.Execute("duncan@oneeyedmen.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$emailParaIndex = $wakelingScope.Paragraphs(1).Index

function Set-ParaPlainText($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $target = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $target.InsertBefore($text)
}

Set-ParaPlainText ($emailParaIndex + 2) "Dr. David Wakeling"
Set-ParaPlainText ($emailParaIndex + 3) "Lecturer at University of Exeter"
Set-ParaPlainText ($emailParaIndex + 4) "D.Wakeling@exeter.ac.uk"
